$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number line and Report Covering the Week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/13/2024  Through  5/19/2024"

# --- Crime data table updates (rows 15-31) ---
# Row 15
$ws.Range("L15").Value = -33.333333333333

# Row 16
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -80
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 39
$ws.Range("J16").Value = 48
$ws.Range("K16").Value = -18.75
$ws.Range("L16").Value = -35
$ws.Range("M16").Value = -17.021276595744
$ws.Range("N16").Value = -86.267605633802

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 53.846153846153
$ws.Range("I17").Value = 96
$ws.Range("J17").Value = 72
$ws.Range("K17").Value = 33.333333333333
$ws.Range("L17").Value = 15.662650602409
$ws.Range("M17").Value = 74.545454545454
$ws.Range("N17").Value = 17.073170731707

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 12
$ws.Range("H18").Value = 9.090909090909
$ws.Range("I18").Value = 53
$ws.Range("J18").Value = 51
$ws.Range("K18").Value = 3.921568627450
$ws.Range("L18").Value = -17.1875
$ws.Range("M18").Value = 140.909090909091
$ws.Range("N18").Value = -60.150375939849

# Row 19
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = -13.157894736842
$ws.Range("I19").Value = 146
$ws.Range("J19").Value = 200
$ws.Range("K19").Value = -27
$ws.Range("L19").Value = -50.170648464163
$ws.Range("M19").Value = 84.810126582278
$ws.Range("N19").Value = -9.316770186335

# Row 20
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 21
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = 23.529411764705
$ws.Range("L20").Value = 10.526315789473
$ws.Range("M20").Value = -8.695652173913
$ws.Range("N20").Value = -85.106382978723

# Row 21
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -20.833333333333
$ws.Range("F21").Value = 79
$ws.Range("G21").Value = 77
$ws.Range("H21").Value = 2.597402597402
$ws.Range("I21").Value = 359
$ws.Range("J21").Value = 393
$ws.Range("K21").Value = -8.651399491094
$ws.Range("L21").Value = -31.878557874762
$ws.Range("M21").Value = 56.768558951965
$ws.Range("N21").Value = -55.514250309789

# Row 22
$ws.Range("C22").Value = "'0"
$ws.Range("D22").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null

# Row 23
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 20
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 82
$ws.Range("J23").Value = 52
$ws.Range("K23").Value = 57.692307692307
$ws.Range("L23").Value = 15.492957746478
$ws.Range("M23").Value = 43.859649122807

# Row 24
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 55
$ws.Range("E24").Value = -49.090909090909
$ws.Range("F24").Value = 106
$ws.Range("G24").Value = 110
$ws.Range("H24").Value = -3.636363636363
$ws.Range("I24").Value = 462
$ws.Range("J24").Value = 426
$ws.Range("K24").Value = 8.450704225352
$ws.Range("L24").Value = -49.782608695652
$ws.Range("M24").Value = 80.46875

# Row 25
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = 69.230769230769
$ws.Range("F25").Value = 72
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = 100
$ws.Range("I25").Value = 298
$ws.Range("J25").Value = 209
$ws.Range("K25").Value = 42.583732057416
$ws.Range("L25").Value = -60.424966799468

# Row 26
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = -64.285714285714
$ws.Range("F26").Value = 32
$ws.Range("G26").Value = 37
$ws.Range("H26").Value = -13.513513513513
$ws.Range("I26").Value = 171
$ws.Range("J26").Value = 153
$ws.Range("K26").Value = 11.764705882352
$ws.Range("L26").Value = -1.156069364161
$ws.Range("M26").Value = 47.413793103448

# Row 27
$ws.Range("D27").Value = "'0"
$ws.Range("C27").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").Value = "'***.*"
$ws.Range("C27").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("L27").Value = -18.181818181818

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D28").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 14
$ws.Range("J28").Value = 23
$ws.Range("K28").Value = -39.130434782608
$ws.Range("L28").Value = -12.5

# Row 29
$ws.Range("L29").Value = -83.333333333333

# Row 30
$ws.Range("L30").Value = -75

# Row 31
$ws.Range("F31").Value = "'0"
$ws.Range("G31").Copy() | Out-Null
$ws.Range("F31").PasteSpecial(-4122) | Out-Null
$ws.Range("I31").Value = 10
$ws.Range("K31").Value = 233.333333333333
$ws.Range("L31").Value = 66.666666666666
